# Generate Report for Handback
#
# The handback-status report was regenerated, which refreshed the
# "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps for the
# 2b34266a-9a16-4b77-96b0-a1a636131231.md file across the Overview,
# zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# Overview sheet: refresh "Latest HO Xliff Generate Date" for the
# 2b34266a-...md row (row 2).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-16 22:47:30"

# zh-cn sheet: refresh Correspond Handoff/Handback Datetime for the
# 2b34266a-...md row (row 2).
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-16 22:47:24"
$wsZhCn.Range("K2").Value = "2016-08-16 22:47:41"

# de-de sheet: refresh Correspond Handoff/Handback Datetime for the
# 2b34266a-...md row (row 2).
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-16 22:47:30"
$wsDeDe.Range("K2").Value = "2016-08-16 22:47:49"
